# BIS-1002: removed "Internal Assignment" column from export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Internal Assignment" header/values in column O (rows 4-7); this
# also drops the now-unused "Internal Assignment" shared string on save.
$ws.Range("O4:O7").ClearContents()

# Scroll/select to match the post-edit view state (selection moved from
# O12 to the cleared O4:O7 range).
$ws.Range("O4:O7").Select()
